$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook stores every data cell as inline/shared TEXT, even when the
# text happens to look like a number (e.g. "37.40", "1.00"). Plain assignment
# of such strings to .Value would let Excel auto-convert them to floating
# point numbers and lose the exact textual formatting (trailing zeros, etc.),
# so for those cells we first force a text ("@") number format.


# Row 2
$ws.Range("D2").Value = "64.290.19"
$ws.Range("E2").Value = "  -1.76%  "

# Row 3
$ws.Range("D3").Value = "3.123.74"
$ws.Range("E3").Value = "  -2.32%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.67"
$ws.Range("E5").Value = "  -0.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.40"
$ws.Range("E6").Value = "  +2.94%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").Value = "3.122.28"
$ws.Range("E9").Value = "  -2.34%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("E10").Value = "  -5.36%  "

# Row 11
$ws.Range("E11").Value = "  -2.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -3.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.40"
$ws.Range("E13").Value = "  -4.93%  "

# Row 14
$ws.Range("E14").Value = "  -5.79%  "

# Row 15
$ws.Range("D15").Value = "3.641.73"
$ws.Range("E15").Value = "  -2.42%  "

# Row 16
$ws.Range("E16").Value = "  -1.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.29"
$ws.Range("E17").Value = "  -1.81%  "

# Row 18
$ws.Range("D18").Value = "64.237.48"

# Row 19
$ws.Range("D19").Value = "3.125.82"
$ws.Range("E19").Value = "  -2.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "479.17"
$ws.Range("E20").Value = "  -1.22%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.59"
$ws.Range("E21").Value = "  -3.25%  "

# Row 22
$ws.Range("E22").Value = "  -6.80%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  -3.86%  "

# Row 24
$ws.Range("E24").Value = "  +1.51%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.02"
$ws.Range("E25").Value = "  -6.61%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.51"
$ws.Range("E26").Value = "  -2.50%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.55"
$ws.Range("E27").Value = "  +7.45%  "

# Row 28
$ws.Range("E28").Value = "  -0.32%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.64"
$ws.Range("E29").Value = "  +2.14%  "

# Row 30
$ws.Range("E30").Value = "  -2.78%  "

# Row 31
$ws.Range("E31").Value = "  -2.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.06%  "

# Row 33
$ws.Range("E33").Value = "  -5.75%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.41"
$ws.Range("E34").Value = "  -4.17%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0850"
$ws.Range("E35").Value = "  -5.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  -2.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.34"
$ws.Range("E37").Value = "  -7.99%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.06"
$ws.Range("E38").Value = "  -4.62%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("E39").Value = "  -5.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.20"
$ws.Range("E40").Value = "  -0.83%  "

# Row 41
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.19"
$ws.Range("E41").Value = "  -2.72%  "

# Row 42
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "449.93"
$ws.Range("E42").Value = "  -5.55%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.293"
$ws.Range("E43").Value = "  -2.91%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0367"
$ws.Range("E44").Value = "  -4.37%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.113"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.23"
$ws.Range("E46").Value = "  +3.93%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.843.07"
$ws.Range("E47").Value = "  -3.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.88"
$ws.Range("E48").Value = "  -0.68%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.84"
$ws.Range("E49").Value = "  +0.97%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("E51").Value = "  -2.82%  "
